$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 24888
$ws.Range("J3").Value = 24888
$ws.Range("L3").Value = 24888
$ws.Range("N3").Value = -25116
$ws.Range("H33").Value = 297.03705
$ws.Range("I33").Value = 313.2
$ws.Range("K33").Value = 313.2
$ws.Range("M33").Value = -84.19999999999999
$ws.Range("H55").Value = 1195.5714
$ws.Range("I55").Value = 2440.2
$ws.Range("J55").Value = 504.1111
$ws.Range("K55").Value = 2440.2
$ws.Range("L55").Value = 504.1111
$ws.Range("M55").Value = -2226.2
$ws.Range("N55").Value = -932.1111000000001
$ws.Range("H102").Value = 24888
$ws.Range("J102").Value = 24888
$ws.Range("L102").Value = 24888
$ws.Range("N102").Value = -31378
$ws.Range("H116").Value = 2749.25
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 2974.3845
$ws.Range("I132").Value = 2946.3044
$ws.Range("K132").Value = 8838.913199999999
$ws.Range("M132").Value = -6308.913199999999
$ws.Range("H138").Value = 5197
$ws.Range("I138").Value = 5197
$ws.Range("K138").Value = 15591
$ws.Range("M138").Value = -10451
$ws.Range("H141").Value = 6241.077
$ws.Range("I141").Value = 6682.1665
$ws.Range("J141").Value = 948
$ws.Range("K141").Value = 20046.4995
$ws.Range("L141").Value = 2844
$ws.Range("M141").Value = -14866.4995
$ws.Range("N141").Value = -13204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1139.6666
$ws.Range("I16").Value = 1139.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1139.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -852.6666
$ws.Range("N16").Value = ""
$ws.Range("H32").Value = 1678.579
$ws.Range("I32").Value = 1699.75
$ws.Range("J32").Value = 1297.5
$ws.Range("K32").Value = 1699.75
$ws.Range("L32").Value = 1297.5
$ws.Range("M32").Value = -1412.75
$ws.Range("N32").Value = -1871.5
$ws.Range("H96").Value = 26081.5
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("H132").Value = 2760.5789
$ws.Range("I132").Value = 2760.5789
$ws.Range("K132").Value = 8281.736699999999
$ws.Range("M132").Value = -5751.736699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2508.4443
$ws.Range("I20").Value = 2927.5
$ws.Range("J20").Value = 2173.2
$ws.Range("K20").Value = 2927.5
$ws.Range("L20").Value = 2173.2
$ws.Range("M20").Value = -2680.5
$ws.Range("N20").Value = -2667.2
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15478
$ws.Range("H80").Value = 4351.5
$ws.Range("I80").Value = 464.36365
$ws.Range("J80").Value = 10459.857
$ws.Range("K80").Value = 464.36365
$ws.Range("L80").Value = 10459.857
$ws.Range("M80").Value = 533.63635
$ws.Range("N80").Value = -12455.857
$ws.Range("H83").Value = 4351.5
$ws.Range("I83").Value = 464.36365
$ws.Range("J83").Value = 10459.857
$ws.Range("K83").Value = 2321.81825
$ws.Range("L83").Value = 52299.285
$ws.Range("M83").Value = 2670.18175
$ws.Range("N83").Value = -62283.285
$ws.Range("H134").Value = 7665.077
$ws.Range("I134").Value = 7470.6665
$ws.Range("K134").Value = 22411.9995
$ws.Range("M134").Value = -19876.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5000572.5
$ws.Range("I22").Value = 470.25
$ws.Range("J22").Value = 10000675
$ws.Range("K22").Value = 470.25
$ws.Range("L22").Value = 10000675
$ws.Range("M22").Value = -120.25
$ws.Range("N22").Value = -10001375
$ws.Range("H132").Value = 690.1053000000001
$ws.Range("I132").Value = 676.125
$ws.Range("K132").Value = 2028.375
$ws.Range("M132").Value = 501.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 270
$ws.Range("I13").Value = 40
$ws.Range("K13").Value = 120
$ws.Range("M13").Value = 48
$ws.Range("H23").Value = 225.83333
$ws.Range("I23").Value = 61
$ws.Range("K23").Value = 183
$ws.Range("M23").Value = 52
$ws.Range("H80").Value = 7084.8335
$ws.Range("I80").Value = 6984.6665
$ws.Range("J80").Value = 7185
$ws.Range("K80").Value = 20953.9995
$ws.Range("L80").Value = 21555
$ws.Range("M80").Value = -20017.9995
$ws.Range("N80").Value = -23427
$ws.Range("H83").Value = 7084.8335
$ws.Range("I83").Value = 6984.6665
$ws.Range("J83").Value = 7185
$ws.Range("K83").Value = 62861.9985
$ws.Range("L83").Value = 64665
$ws.Range("M83").Value = -58181.9985
$ws.Range("N83").Value = -74025
$ws.Range("H86").Value = 943.13336
$ws.Range("I86").Value = 474
$ws.Range("J86").Value = 2819.6667
$ws.Range("K86").Value = 1422
$ws.Range("L86").Value = 8459.000100000001
$ws.Range("M86").Value = -236
$ws.Range("N86").Value = -10831.0001
$ws.Range("H89").Value = 943.13336
$ws.Range("I89").Value = 474
$ws.Range("J89").Value = 2819.6667
$ws.Range("K89").Value = 4266
$ws.Range("L89").Value = 25377.0003
$ws.Range("M89").Value = 1662
$ws.Range("N89").Value = -37233.0003
$ws.Range("H122").Value = 647.6923
$ws.Range("J122").Value = 721.5
$ws.Range("L122").Value = 6493.5
$ws.Range("N122").Value = -11393.5
$ws.Range("H131").Value = 636
$ws.Range("I131").Value = 636
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1908
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 3132
$ws.Range("N131").Value = ""
$ws.Range("H137").Value = 6873.5454
$ws.Range("I137").Value = 7082.5
$ws.Range("J137").Value = 6754.143
$ws.Range("K137").Value = 21247.5
$ws.Range("L137").Value = 20262.429
$ws.Range("M137").Value = -16147.5
$ws.Range("N137").Value = -30462.429
$ws.Range("H139").Value = 4893.0586
$ws.Range("I139").Value = 4812.1333
$ws.Range("K139").Value = 14436.3999
$ws.Range("M139").Value = -9296.3999
$ws.Range("H140").Value = 772561.4
$ws.Range("I140").Value = 772561.4
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 2317684.2
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2312504.2
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 7803.6665
$ws.Range("I141").Value = 7803.6665
$ws.Range("K141").Value = 23410.9995
$ws.Range("M141").Value = -18230.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1472.25
$ws.Range("I31").Value = 1129.6666
$ws.Range("K31").Value = 1129.6666
$ws.Range("M31").Value = -837.6666
$ws.Range("H37").Value = 1472.25
$ws.Range("I37").Value = 1129.6666
$ws.Range("K37").Value = 1129.6666
$ws.Range("M37").Value = -852.6666
$ws.Range("H46").Value = 21499.533
$ws.Range("I46").Value = 5868.3335
$ws.Range("J46").Value = 25407.334
$ws.Range("K46").Value = 5868.3335
$ws.Range("L46").Value = 25407.334
$ws.Range("M46").Value = -5712.3335
$ws.Range("N46").Value = -25719.334
$ws.Range("H122").Value = 2718
$ws.Range("I122").Value = 2269.3333
$ws.Range("J122").Value = 3166.6667
$ws.Range("K122").Value = 6807.999899999999
$ws.Range("L122").Value = 9500.000100000001
$ws.Range("M122").Value = -4357.999899999999
$ws.Range("N122").Value = -14400.0001
$ws.Range("H132").Value = 2590.75
$ws.Range("I132").Value = 2496.8
$ws.Range("K132").Value = 7490.400000000001
$ws.Range("M132").Value = -4960.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H68").Value = 3004.2856
$ws.Range("I68").Value = 2764.3333
$ws.Range("J68").Value = 4444
$ws.Range("K68").Value = 2764.3333
$ws.Range("L68").Value = 4444
$ws.Range("M68").Value = -2015.3333
$ws.Range("N68").Value = -5942
$ws.Range("H71").Value = 3004.2856
$ws.Range("I71").Value = 2764.3333
$ws.Range("J71").Value = 4444
$ws.Range("K71").Value = 13821.6665
$ws.Range("L71").Value = 22220
$ws.Range("M71").Value = -10077.6665
$ws.Range("N71").Value = -29708
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9999.5
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 9999
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -9375
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 9999.5
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 49995
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -46875
$ws.Range("N65").Value = -56240
$ws.Range("H81").Value = 3314.0527
$ws.Range("I81").Value = 3292.2942
$ws.Range("J81").Value = 3499
$ws.Range("K81").Value = 6584.5884
$ws.Range("L81").Value = 6998
$ws.Range("M81").Value = -5523.5884
$ws.Range("N81").Value = -9120
$ws.Range("H84").Value = 3314.0527
$ws.Range("I84").Value = 3292.2942
$ws.Range("J84").Value = 3499
$ws.Range("K84").Value = 32922.942
$ws.Range("L84").Value = 34990
$ws.Range("M84").Value = -27618.942
$ws.Range("N84").Value = -45598
$ws.Range("H112").Value = 30247.25
$ws.Range("J112").Value = 30247.25
$ws.Range("L112").Value = 30247.25
$ws.Range("N112").Value = -33201.25
$ws.Range("H122").Value = 4698.5
$ws.Range("I122").Value = 4197.9287
$ws.Range("J122").Value = 8202.5
$ws.Range("K122").Value = 12593.7861
$ws.Range("L122").Value = 24607.5
$ws.Range("M122").Value = -10143.7861
$ws.Range("N122").Value = -29507.5
$ws.Range("H132").Value = 3043.2856
$ws.Range("I132").Value = 1742.8462
$ws.Range("K132").Value = 5228.5386
$ws.Range("M132").Value = -2698.5386
